$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "24-03-2025"
$ws.Range("B5").Value = "Delhi Capitals vs Lucknow Super Giants"
$ws.Range("C5").Value = "Delhi Capitals"
$ws.Range("D5").Value = "Delhi Capitals"
